$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 378, pushing existing rows 378:467 down to 379:468
$ws.Rows("378").Insert()

# Populate the new row 378 with the new weekly record
$ws.Range("A378").Value = 9
$ws.Range("B378").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C378").Value = 'Metropolitana'
$ws.Range("D378").Value = 45258
$ws.Range("D378").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E378").Value = 13
$ws.Range("F378").Value = 100112021
$ws.Range("G378").Value = 'Ají'
$ws.Range("H378").Value = 'Americana (o)'
$ws.Range("I378").Value = 'Primera'
$ws.Range("J378").Value = 45
$ws.Range("K378").Value = 44000
$ws.Range("L378").Value = 45000
$ws.Range("M378").Value = 44444
$ws.Range("N378").Value = '$/caja 25 kilos'
$ws.Range("O378").Value = 'Provincia de Limarí'
$ws.Range("P378").Value = 1778
$ws.Range("Q378").Value = 25
$ws.Range("R378").Value = 'Hortaliza'
